# remove accel filter in detect_kneepoint
# Add a new "CUTOFF_FREQ" parameter row to the params sheet (used by the
# InputHandler), inserted right after the RE-SAMPLE row and before PRE_TIME.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")
$ws.Activate()

# Insert a new row at row 3, pushing PRE_TIME and everything below it down.
$ws.Rows.Item(3).Insert()

# Populate the new row: Parameter, Value, Type, Unit, Description, Class
$ws.Cells.Item(3, 1).Value = "CUTOFF_FREQ"
$ws.Cells.Item(3, 2).Value = 10
$ws.Cells.Item(3, 3).Value = "int"
$ws.Cells.Item(3, 6).Value = "InputHandler"

# Match the view/selection state recorded in the saved workbook.
$ws.Range("E20").Select()

# Match the page setup recorded in the saved workbook.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
